# Masamune_Profits data refresh: updated currentAveragePrice / LevePrice
# feed values (and their derived profit columns) across the 8 crafting-job
# sheets, as pulled by the scheduled Sheets-update runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 77008970
$ws.Range("I96").Value = 6926.4
$ws.Range("J96").Value = 125135250
$ws.Range("K96").Value = 20779.2
$ws.Range("L96").Value = 375405750
$ws.Range("M96").Value = -19406.2
$ws.Range("N96").Value = -375408496

$ws.Range("H129").Value = 1085.28
$ws.Range("I129").Value = 780.6923
$ws.Range("J129").Value = 1130.7931
$ws.Range("K129").Value = 2342.0769
$ws.Range("L129").Value = 3392.379300000001
$ws.Range("M129").Value = 2657.9231
$ws.Range("N129").Value = -13392.3793

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 100510
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 200020
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 200020
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -203264

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 18625
$ws.Range("J93").Value = 18625
$ws.Range("L93").Value = 18625
$ws.Range("N93").Value = -22369

$ws.Range("H105").Value = 3602
$ws.Range("I105").Value = 4150.25
$ws.Range("J105").Value = 2505.5
$ws.Range("K105").Value = 4150.25
$ws.Range("L105").Value = 2505.5
$ws.Range("M105").Value = -2403.25
$ws.Range("N105").Value = -5999.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 68571.42999999999
$ws.Range("J68").Value = 68571.42999999999
$ws.Range("L68").Value = 68571.42999999999
$ws.Range("N68").Value = -70069.42999999999

$ws.Range("H71").Value = 68571.42999999999
$ws.Range("J71").Value = 68571.42999999999
$ws.Range("L71").Value = 205714.29
$ws.Range("N71").Value = -213202.29

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H80").Value = 36763.5
$ws.Range("J80").Value = 36763.5
$ws.Range("L80").Value = 36763.5
$ws.Range("N80").Value = -39009.5

$ws.Range("H83").Value = 36763.5
$ws.Range("J83").Value = 36763.5
$ws.Range("L83").Value = 110290.5
$ws.Range("N83").Value = -121522.5

$ws.Range("H99").Value = 2457.4
$ws.Range("I99").Value = 2412
$ws.Range("J99").Value = 2468.75
$ws.Range("K99").Value = 2412
$ws.Range("L99").Value = 2468.75
$ws.Range("M99").Value = -914
$ws.Range("N99").Value = -5464.75

$ws.Range("H126").Value = 2457.4
$ws.Range("I126").Value = 2412
$ws.Range("J126").Value = 2468.75
$ws.Range("K126").Value = 7236
$ws.Range("L126").Value = 7406.25
$ws.Range("M126").Value = -4766
$ws.Range("N126").Value = -12346.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3990.1177
$ws.Range("I3").Value = 1807.6666
$ws.Range("J3").Value = 6445.375
$ws.Range("K3").Value = 5422.9998
$ws.Range("L3").Value = 19336.125
$ws.Range("M3").Value = -5310.9998
$ws.Range("N3").Value = -19560.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H80").Value = 459635.9
$ws.Range("I80").Value = 504799.5
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 504799.5
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -503801.5
$ws.Range("N80").Value = -9996

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H83").Value = 459635.9
$ws.Range("I83").Value = 504799.5
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 2523997.5
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -2519005.5
$ws.Range("N83").Value = -49984

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 40000
$ws.Range("J18").Value = 40000
$ws.Range("L18").Value = 40000
$ws.Range("N18").Value = -40344

$ws.Range("H22").Value = 1013.1818
$ws.Range("I22").Value = 986.9375
$ws.Range("J22").Value = 1083.1666
$ws.Range("K22").Value = 986.9375
$ws.Range("L22").Value = 1083.1666
$ws.Range("M22").Value = -691.9375
$ws.Range("N22").Value = -1673.1666

$ws.Range("H27").Value = 1013.1818
$ws.Range("I27").Value = 986.9375
$ws.Range("J27").Value = 1083.1666
$ws.Range("K27").Value = 986.9375
$ws.Range("L27").Value = 1083.1666
$ws.Range("M27").Value = -879.9375
$ws.Range("N27").Value = -1297.1666

$ws.Range("H46").Value = 3848.8572
$ws.Range("I46").Value = 1293.3334
$ws.Range("J46").Value = 4545.8184
$ws.Range("K46").Value = 1293.3334
$ws.Range("L46").Value = 4545.8184
$ws.Range("M46").Value = -1105.3334
$ws.Range("N46").Value = -4921.8184

$ws.Range("H68").Value = 3146.7334
$ws.Range("I68").Value = 2933.5
$ws.Range("J68").Value = 3288.889
$ws.Range("K68").Value = 2933.5
$ws.Range("L68").Value = 3288.889
$ws.Range("M68").Value = -2184.5
$ws.Range("N68").Value = -4786.889

$ws.Range("H71").Value = 3146.7334
$ws.Range("I71").Value = 2933.5
$ws.Range("J71").Value = 3288.889
$ws.Range("K71").Value = 14667.5
$ws.Range("L71").Value = 16444.445
$ws.Range("M71").Value = -10923.5
$ws.Range("N71").Value = -23932.445

$ws.Range("H82").Value = 8335120.5
$ws.Range("I82").Value = 1978
$ws.Range("J82").Value = 16668263
$ws.Range("K82").Value = 1978
$ws.Range("L82").Value = 16668263
$ws.Range("M82").Value = -1617
$ws.Range("N82").Value = -16668985

$ws.Range("H85").Value = 8335120.5
$ws.Range("I85").Value = 1978
$ws.Range("J85").Value = 16668263
$ws.Range("K85").Value = 1978
$ws.Range("L85").Value = 16668263
$ws.Range("M85").Value = -730
$ws.Range("N85").Value = -16670759

$ws.Range("H99").Value = 27998.666
$ws.Range("J99").Value = 29998.4
$ws.Range("L99").Value = 29998.4
$ws.Range("N99").Value = -35988.4

$ws.Range("H122").Value = 60699.65
$ws.Range("I122").Value = 101619.4
$ws.Range("J122").Value = 2242.8572
$ws.Range("K122").Value = 304858.2
$ws.Range("L122").Value = 6728.571599999999
$ws.Range("M122").Value = -302408.2
$ws.Range("N122").Value = -11628.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2893
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2893
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H93").Value = 35100.715
$ws.Range("J93").Value = 35100.715
$ws.Range("L93").Value = 35100.715
$ws.Range("N93").Value = -40092.715
